# Update column F (CE) values per data refresh on 03.03.2019
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(30, 6).Value = -3.262744396099325
$ws.Cells.Item(34, 6).Value = -3.261953322486035
$ws.Cells.Item(39, 6).Value = -3.244558414039905
$ws.Cells.Item(47, 6).Value = -3.181603642213611
$ws.Cells.Item(49, 6).Value = -3.163001107696576
$ws.Cells.Item(52, 6).Value = -3.136528714473825
$ws.Cells.Item(53, 6).Value = -3.124533961246978
$ws.Cells.Item(79, 6).Value = -3.329428529429942
$ws.Cells.Item(80, 6).Value = -3.33683976362203
$ws.Cells.Item(81, 6).Value = -3.343757573777839
$ws.Cells.Item(82, 6).Value = -3.349121933240415
$ws.Cells.Item(83, 6).Value = -3.35065255953292
$ws.Cells.Item(85, 6).Value = -3.362294303171793
$ws.Cells.Item(87, 6).Value = -3.364813611452816
$ws.Cells.Item(89, 6).Value = -3.372126827939522
$ws.Cells.Item(92, 6).Value = -3.379625888497094
$ws.Cells.Item(93, 6).Value = -3.387558032591784
$ws.Cells.Item(94, 6).Value = -3.385993764268941
$ws.Cells.Item(97, 6).Value = -3.410581453314329
$ws.Cells.Item(98, 6).Value = -3.408845776635839
$ws.Cells.Item(100, 6).Value = -3.407180173167276
$ws.Cells.Item(101, 6).Value = -3.413052841253143
$ws.Cells.Item(105, 6).Value = -3.422356476545646
$ws.Cells.Item(106, 6).Value = -3.430806461066061
$ws.Cells.Item(108, 6).Value = -3.425103157036734
$ws.Cells.Item(110, 6).Value = -3.432427179511175
$ws.Cells.Item(115, 6).Value = -3.431489785727089
$ws.Cells.Item(119, 6).Value = -3.430220392175171
$ws.Cells.Item(121, 6).Value = -3.425803169821137
$ws.Cells.Item(122, 6).Value = -3.421453714532021
$ws.Cells.Item(124, 6).Value = -3.420803310112217
$ws.Cells.Item(125, 6).Value = -3.416382471003694
$ws.Cells.Item(126, 6).Value = -3.420889410914565
$ws.Cells.Item(127, 6).Value = -3.408569648477553
$ws.Cells.Item(130, 6).Value = -3.403902446737506
$ws.Cells.Item(131, 6).Value = -3.406454779000018
$ws.Cells.Item(134, 6).Value = -3.398794652215288
$ws.Cells.Item(136, 6).Value = -3.396402503187317
$ws.Cells.Item(137, 6).Value = -3.394684730126925
$ws.Cells.Item(138, 6).Value = -3.392369982454565
$ws.Cells.Item(141, 6).Value = -3.385229504827763
$ws.Cells.Item(142, 6).Value = -3.379516993550525
$ws.Cells.Item(144, 6).Value = -3.373368181227372
$ws.Cells.Item(147, 6).Value = -3.360110207630812
$ws.Cells.Item(150, 6).Value = -3.355747737789049
$ws.Cells.Item(151, 6).Value = -3.352094058307857
$ws.Cells.Item(152, 6).Value = -3.346526254866567
$ws.Cells.Item(155, 6).Value = -3.332571949820299
$ws.Cells.Item(158, 6).Value = -3.327985667701459
$ws.Cells.Item(160, 6).Value = -3.320443749766544
$ws.Cells.Item(161, 6).Value = -3.315823554375775
$ws.Cells.Item(167, 6).Value = -3.290622586515872
$ws.Cells.Item(168, 6).Value = -3.285792181869756
$ws.Cells.Item(169, 6).Value = -3.282186816589732
$ws.Cells.Item(170, 6).Value = -3.276957631474528
$ws.Cells.Item(172, 6).Value = -3.268925143700635
$ws.Cells.Item(175, 6).Value = -3.254511213314867
$ws.Cells.Item(176, 6).Value = -3.2515119486168
$ws.Cells.Item(177, 6).Value = -3.24675010450867
$ws.Cells.Item(178, 6).Value = -3.242129253773831
$ws.Cells.Item(180, 6).Value = -3.233320556388605
$ws.Cells.Item(184, 6).Value = -3.215826186402178
$ws.Cells.Item(186, 6).Value = -3.205838286710128
$ws.Cells.Item(188, 6).Value = -3.196174925523536
$ws.Cells.Item(190, 6).Value = -3.185818856042483
$ws.Cells.Item(191, 6).Value = -3.183218800371353
$ws.Cells.Item(192, 6).Value = -3.175993698592579
$ws.Cells.Item(193, 6).Value = -3.171067706065466
$ws.Cells.Item(194, 6).Value = -3.165705596181808
$ws.Cells.Item(195, 6).Value = -3.15987713428439
$ws.Cells.Item(196, 6).Value = -3.154090072665406
$ws.Cells.Item(198, 6).Value = -3.143827441834141
$ws.Cells.Item(200, 6).Value = -3.132074977726777
$ws.Cells.Item(201, 6).Value = -3.127115903630184
$ws.Cells.Item(205, 6).Value = -3.103700145850151
$ws.Cells.Item(206, 6).Value = -3.097823913796469
$ws.Cells.Item(207, 6).Value = -3.091947681742789
$ws.Cells.Item(221, 6).Value = -3.450747757291689
$ws.Cells.Item(222, 6).Value = -3.498223947086537
$ws.Cells.Item(223, 6).Value = -3.512746184054146
$ws.Cells.Item(224, 6).Value = -3.492254115884029
$ws.Cells.Item(225, 6).Value = -3.476079015660526
$ws.Cells.Item(226, 6).Value = -3.41683062528953
$ws.Cells.Item(227, 6).Value = -3.357931575278279
$ws.Cells.Item(229, 6).Value = -3.210910031589647
$ws.Cells.Item(232, 6).Value = -3.524536394082353
$ws.Cells.Item(233, 6).Value = -3.550261960062457
$ws.Cells.Item(234, 6).Value = -3.546950400261686
$ws.Cells.Item(235, 6).Value = -3.556012789047096
$ws.Cells.Item(236, 6).Value = -3.51657924830859
$ws.Cells.Item(237, 6).Value = -3.480575700030062
$ws.Cells.Item(238, 6).Value = -3.419422554187102
$ws.Cells.Item(243, 6).Value = -3.570536380760595
$ws.Cells.Item(244, 6).Value = -3.590321433278569
$ws.Cells.Item(245, 6).Value = -3.599820970244835
$ws.Cells.Item(246, 6).Value = -3.601076490031761
$ws.Cells.Item(247, 6).Value = -3.561473773757928
$ws.Cells.Item(248, 6).Value = -3.513661946763981
$ws.Cells.Item(250, 6).Value = -3.396944368536094
$ws.Cells.Item(251, 6).Value = -3.322503860871483
$ws.Cells.Item(254, 6).Value = -3.600382309580585
$ws.Cells.Item(255, 6).Value = -3.620815515683325
$ws.Cells.Item(257, 6).Value = -3.621767690297197
$ws.Cells.Item(258, 6).Value = -3.58054310001445
$ws.Cells.Item(259, 6).Value = -3.540984566425637
$ws.Cells.Item(260, 6).Value = -3.490152088682049
$ws.Cells.Item(261, 6).Value = -3.427879325927404
$ws.Cells.Item(265, 6).Value = -3.625893799373868
$ws.Cells.Item(266, 6).Value = -3.640413019611778
$ws.Cells.Item(267, 6).Value = -3.648650445159416
$ws.Cells.Item(268, 6).Value = -3.636534180483225
$ws.Cells.Item(269, 6).Value = -3.607688031948528
$ws.Cells.Item(270, 6).Value = -3.568590549616089
$ws.Cells.Item(271, 6).Value = -3.513818264561325
$ws.Cells.Item(272, 6).Value = -3.450125577162325
$ws.Cells.Item(273, 6).Value = -3.381310845999874
$ws.Cells.Item(276, 6).Value = -3.646660509348909
$ws.Cells.Item(277, 6).Value = -3.66246163278088
$ws.Cells.Item(278, 6).Value = -3.660169587957068
$ws.Cells.Item(279, 6).Value = -3.645764315762032
$ws.Cells.Item(280, 6).Value = -3.622195376379684
$ws.Cells.Item(281, 6).Value = -3.582411526465713
$ws.Cells.Item(282, 6).Value = -3.533432120910382
$ws.Cells.Item(284, 6).Value = -3.400859203289268
$ws.Cells.Item(287, 6).Value = -3.662505098397937
$ws.Cells.Item(288, 6).Value = -3.673389180865456
$ws.Cells.Item(289, 6).Value = -3.672734372040193
$ws.Cells.Item(290, 6).Value = -3.658041505470925
$ws.Cells.Item(291, 6).Value = -3.631747326963702
$ws.Cells.Item(292, 6).Value = -3.596923147439703
$ws.Cells.Item(293, 6).Value = -3.546416154189119
$ws.Cells.Item(294, 6).Value = -3.484721573868963

Write-Output "Updated 132 cells in column F"
